$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new Mac-Addresses (user rows) appended to the table: "John Doe" and
# "Jane Smith". Write the brand-new text values in the same order they were
# added to the workbook's shared-string table (John Doe's name/e-mail first,
# then Jane Smith's) even though Jane Smith physically lands on the earlier
# row (31) and John Doe on row 32.
$ws.Range("C32").Value = "John Doe"
$ws.Range("D32").Value = "john.doe@xyz.com"
$ws.Range("C31").Value = "Jane Smith"
$ws.Range("D31").Value = "jane.smith@xyz.com"

# Row 31
$ws.Range("A31").Value = 110030
$ws.Range("B31").Value = 9317596768
$ws.Range("E31").Value = 818876432
$ws.Range("F31").Value = "ACT"
$ws.Range("G31").Value = "eng"
$ws.Range("H31").Value = "PWD"
$ws.Range("I31").Value = $true
$ws.Range("J31").Value = "superadmin"
$ws.Range("K31").Value = "now()"

# Row 32
$ws.Range("A32").Value = 110031
$ws.Range("B32").Value = 9317596767
$ws.Range("E32").Value = 818876431
$ws.Range("F32").Value = "ACT"
$ws.Range("G32").Value = "eng"
$ws.Range("H32").Value = "PWD"
$ws.Range("I32").Value = $true
$ws.Range("J32").Value = "superadmin"
$ws.Range("K32").Value = "now()"

# Match the left-aligned formatting used by the is_active column for every
# other row in the table.
$ws.Range("I31").HorizontalAlignment = -4131
$ws.Range("I32").HorizontalAlignment = -4131

# Select E28, matching the saved workbook's final selection.
$ws.Range("E28").Select() | Out-Null
